$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,(2, "ECs", "Mmp13", "Lrp1", "ECs", 2, 0.6666666666666666, 0.319463, 0.9583889999999999, 0.004642769429023782, 0.004642769429023782, 3, 1, 23.63579766666667, 70.907393, 0.06827844587621175, 0.06827844587621175, 7.550762829986333, 67.95686546987699, 0.0003170010811753309, 0.0003170010811753309)
    ,(3, "ECs", "Mmp13", "Lrp1", "FAPs", 2, 0.6666666666666666, 0.319463, 0.9583889999999999, 0.004642769429023782, 0.004642769429023782, 3, 1, 181.2883913333334, 543.865174, 0.5237009467675041, 0.523700946767504, 57.91493336052067, 521.234400244686, 0.002431422745602979, 0.002431422745602979)
    ,(4, "ECs", "Mmp13", "Lrp1", "M2", 2, 0.6666666666666666, 0.319463, 0.9583889999999999, 0.004642769429023782, 0.004642769429023782, 3, 1, 111.1005463333333, 333.301639, 0.3209442197221123, 0.3209442197221123, 35.49251383328566, 319.4326244995709, 0.001490070011747715, 0.001490070011747715)
    ,(5, "ECs", "Mmp13", "Lrp1", "sCs", 2, 0.6666666666666666, 0.319463, 0.9583889999999999, 0.004642769429023782, 0.004642769429023782, 3, 1, 30.14303933333333, 90.42911799999999, 0.08707638763417187, 0.08707638763417187, 9.629585774544665, 86.66627197090199, 0.0004042755904977577, 0.0004042755904977577)
    ,(6, "FAPs", "Mmp13", "Lrp1", "ECs", 3, 1, 62.87584200000001, 188.627526, 0.913777298346693, 0.913777298346693, 3, 1, 23.63579766666667, 70.907393, 0.06827844587621175, 0.06827844587621175, 1486.120679633302, 13375.08611669972, 0.06239129380807567, 0.06239129380807567)
    ,(7, "FAPs", "Mmp13", "Lrp1", "FAPs", 3, 1, 62.87584200000001, 188.627526, 0.913777298346693, 0.913777298346693, 3, 1, 181.2883913333334, 543.865174, 0.5237009467675041, 0.523700946767504, 11398.66024990884, 102587.9422491795, 0.4785460362788152, 0.4785460362788151)
    ,(8, "FAPs", "Mmp13", "Lrp1", "M2", 3, 1, 62.87584200000001, 188.627526, 0.913777298346693, 0.913777298346693, 3, 1, 111.1005463333333, 333.301639, 0.3209442197221123, 0.3209442197221123, 6985.540397368346, 62869.86357631511, 0.2932715420176592, 0.2932715420176592)
    ,(9, "FAPs", "Mmp13", "Lrp1", "sCs", 3, 1, 62.87584200000001, 188.627526, 0.913777298346693, 0.913777298346693, 3, 1, 30.14303933333333, 90.42911799999999, 0.08707638763417187, 0.08707638763417187, 1895.268978522452, 17057.42080670207, 0.07956842624214296, 0.07956842624214296)
    ,(10, "M2", "Mmp13", "Lrp1", "ECs", 3, 1, 2.989864666666667, 8.969594, 0.04345183095168575, 0.04345183095168574, 3, 1, 23.63579766666667, 70.907393, 0.06827844587621175, 0.06827844587621175, 70.66783631204912, 636.010526808442, 0.002966823487856978, 0.002966823487856977)
    ,(11, "M2", "Mmp13", "Lrp1", "FAPs", 3, 1, 2.989864666666667, 8.969594, 0.04345183095168575, 0.04345183095168574, 3, 1, 181.2883913333334, 543.865174, 0.5237009467675041, 0.523700946767504, 542.027755724373, 4878.249801519357, 0.02275576500817936, 0.02275576500817936)
    ,(12, "M2", "Mmp13", "Lrp1", "M2", 3, 1, 2.989864666666667, 8.969594, 0.04345183095168575, 0.04345183095168574, 3, 1, 111.1005463333333, 333.301639, 0.3209442197221123, 0.3209442197221123, 332.1755979293962, 2989.580381364566, 0.01394561398028591, 0.01394561398028591)
    ,(13, "M2", "Mmp13", "Lrp1", "sCs", 3, 1, 2.989864666666667, 8.969594, 0.04345183095168575, 0.04345183095168574, 3, 1, 30.14303933333333, 90.42911799999999, 0.08707638763417187, 0.08707638763417187, 90.12360824867689, 811.112474238092, 0.003783628475363495, 0.003783628475363495)
    ,(14, "sCs", "Mmp13", "Lrp1", "ECs", 3, 1, 2.623545666666667, 7.870637, 0.03812810127259751, 0.03812810127259751, 3, 1, 23.63579766666667, 70.907393, 0.06827844587621175, 0.06827844587621175, 62.00959454659344, 558.086350919341, 0.002603327499103769, 0.002603327499103769)
    ,(15, "sCs", "Mmp13", "Lrp1", "FAPs", 3, 1, 2.623545666666667, 7.870637, 0.03812810127259751, 0.03812810127259751, 3, 1, 181.2883913333334, 543.865174, 0.5237009467675041, 0.523700946767504, 475.6183734995376, 4280.565361495838, 0.01996772273490659, 0.01996772273490659)
    ,(16, "sCs", "Mmp13", "Lrp1", "M2", 3, 1, 2.623545666666667, 7.870637, 0.03812810127259751, 0.03812810127259751, 3, 1, 111.1005463333333, 333.301639, 0.3209442197221123, 0.3209442197221123, 291.4773568971159, 2623.296212074043, 0.01223699371241948, 0.01223699371241948)
    ,(17, "sCs", "Mmp13", "Lrp1", "sCs", 3, 1, 2.623545666666667, 7.870637, 0.03812810127259751, 0.03812810127259751, 3, 1, 30.14303933333333, 90.42911799999999, 0.08707638763417187, 0.08707638763417187, 79.08164022312955, 711.734762008166, 0.003320057326167663, 0.003320057326167663)
)

foreach ($r in $data) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 1).Value = $r[1]
    $ws.Cells.Item($rowNum, 2).Value = $r[2]
    $ws.Cells.Item($rowNum, 3).Value = $r[3]
    $ws.Cells.Item($rowNum, 4).Value = $r[4]
    for ($col = 5; $col -le 20; $col++) {
        $ws.Cells.Item($rowNum, $col).Value = $r[$col]
    }
}